$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = (Get-Date -Year 2022 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100114002
$ws.Cells.Item(13, 7).Value = "Camote"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 40
$ws.Cells.Item(13, 11).Value = 24000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 24000
$ws.Cells.Item(13, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 1200
$ws.Cells.Item(13, 17).Value = 20
$ws.Cells.Item(13, 18).Value = "Hortaliza"
